# Update the carjacking-by-neighborhood-by-month workbook with the
# newly-added 2021-11-19 data pull (through November 11 instead of
# November 10).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet tab and refresh the "through" label used in the
# header row / shared string table.
$ws.Name = "Through 2021-11-11"
$ws.Range("B1").Value = "November 2021 (through November 11)"

# Update the current-month (column B) and other historical counts that
# changed for specific neighborhoods (rows).
$ws.Range("M2").Value = 8      # North Lawndale, November 2020
$ws.Range("AT2").Value = 2     # North Lawndale, November 2017

$ws.Range("M3").Value = 6      # Garfield Park, November 2020

$ws.Range("AT4").Value = 7     # Austin, November 2017

$ws.Range("X6").Value = 1      # West Town, November 2019

$ws.Range("B7").Value = 4      # Englewood, current month
$ws.Range("M7").Value = 3      # Englewood, November 2020

$ws.Range("B8").Value = 2      # South Shore, current month
$ws.Range("M8").Value = 4      # South Shore, November 2020

$ws.Range("AT9").Value = 7     # Auburn Gresham, November 2017
$ws.Range("BP9").Value = 3     # Auburn Gresham, November 2015

$ws.Range("B14").Value = 2     # Wicker Park, current month

$ws.Range("AT17").Value = 3    # West Loop, November 2017

$ws.Range("X22").Value = 1     # Little Village, November 2019

$ws.Range("BE42").Value = 1    # Bucktown, November 2016

$ws.Range("B44").Value = 2     # United Center, current month
$ws.Range("BE44").Value = 1    # United Center, November 2016

$ws.Range("BE57").Value = 1    # Lincoln Square, November 2016

$ws.Range("B58").Value = 1     # Fuller Park, current month

$ws.Range("B72").Value = 1     # Gage Park, current month

$ws.Range("M88").Value = 2     # Old Town, November 2020

$ws.Range("AT96").Value = 1    # Ukrainian Village, November 2017

$ws.Range("M98").Value = 2     # Woodlawn, November 2020
